# Shift the Timestamp column (A) forward by 22 days for every data row
# (rows 2-195) and refresh the dependent "Lookup" text in column E, which
# is cached literal text built from the timestamp's date (dd.MM.yyyy) and
# the Quarter number in column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 195; $r++) {
    $timestampCell = $ws.Cells.Item($r, 1)
    $oldTimestamp = $timestampCell.Value()
    $newTimestamp = $oldTimestamp.AddDays(22)
    $timestampCell.Value = $newTimestamp

    $quarter = $ws.Cells.Item($r, 4).Value()
    $lookupCell = $ws.Cells.Item($r, 5)
    $lookupCell.Value = $newTimestamp.ToString("dd.MM.yyyy") + $quarter
}
